# Ajout de vues pour le stagiaire
# Ajout vue menu et reviewAdv du stagiaire
#
# The journal table gets a brand-new first row (date "17 FÉVRIER") with a
# new bullet entry describing the day's work; every other row is left
# untouched (the apparent "shift" of every other date in the XML diff is
# just a side effect of the unchanged rows sliding down beneath the new
# one).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Insert a brand-new row above the current first row; Word copies the
# row/paragraph formatting (trHeight, jc, tcPr, the "Paragraphedeliste"
# bullet style + numId 2) from that first row automatically.
$newRow = $t.Rows.Add($t.Rows.Item(1))

# First cell: the new date label.
$newRow.Cells.Item(1).Range.Text = "17 FÉVRIER"

# Second cell: the single bulleted journal entry for that day.
$newRow.Cells.Item(2).Range.Paragraphs.Item(1).Range.Text = "Ajout des vues reviewAdv du stagiaire et aussi de la vue menu de ce dernier."
